# Update cryptocurrency price/volume figures (refreshed data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.930.15'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '2.753.17'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.03%  '
$cellStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.22'
$ws.Range("D5").Style = $cellStyle
$ws.Range("E5").Value = '  -0.84%  '
$cellStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.57'
$ws.Range("D6").Style = $cellStyle
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("E10").Value = '  +1.26%  '
$cellStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.68'
$ws.Range("D11").Style = $cellStyle
$ws.Range("E11").Value = '  -15.70%  '
$ws.Range("E12").Value = '  -3.07%  '
$ws.Range("D13").Value = '3.238.13'
$ws.Range("E13").Value = '  -0.10%  '
$cellStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.53'
$ws.Range("D14").Style = $cellStyle
$ws.Range("E14").Value = '  -3.88%  '
$ws.Range("D15").Value = '63.546.46'
$ws.Range("E15").Value = '  -0.70%  '
$cellStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000151'
$ws.Range("D16").Style = $cellStyle
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '2.755.11'
$ws.Range("E17").Value = '  -0.44%  '
$cellStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.11'
$ws.Range("D18").Style = $cellStyle
$ws.Range("E18").Value = '  -0.16%  '
$cellStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.82'
$ws.Range("D19").Style = $cellStyle
$ws.Range("E19").Value = '  -2.70%  '
$cellStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '355.21'
$ws.Range("D20").Style = $cellStyle
$ws.Range("E20").Value = '  -2.18%  '
$ws.Range("E21").Value = '  -4.15%  '
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  -0.52%  '
$cellStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.17'
$ws.Range("D24").Style = $cellStyle
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("E26").Value = '  -0.09%  '
$cellStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.47'
$ws.Range("D27").Style = $cellStyle
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("D28").Value = '0.0₃0912'
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  -5.05%  '
$cellStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.08'
$ws.Range("D30").Style = $cellStyle
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("E31").Value = '  -0.99%  '
$cellStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '168.59'
$ws.Range("D32").Style = $cellStyle
$ws.Range("E32").Value = '  -4.02%  '
$ws.Range("E33").Value = '  -2.68%  '
$cellStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.90'
$ws.Range("D34").Style = $cellStyle
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("E37").Value = '  -1.66%  '
$cellStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.986'
$ws.Range("D38").Style = $cellStyle
$ws.Range("E38").Value = '  -2.72%  '
$cellStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.17'
$ws.Range("D39").Style = $cellStyle
$ws.Range("E39").Value = '  +6.29%  '
$cellStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '331.57'
$ws.Range("D40").Style = $cellStyle
$ws.Range("E40").Value = '  -2.62%  '
$cellStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.15'
$ws.Range("D41").Style = $cellStyle
$ws.Range("E41").Value = '  -4.06%  '
$ws.Range("E42").Value = '  -0.98%  '
$cellStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.47'
$ws.Range("D43").Style = $cellStyle
$ws.Range("E43").Value = '  -2.36%  '
$cellStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0588'
$ws.Range("D44").Style = $cellStyle
$ws.Range("E44").Value = '  -2.12%  '
$cellStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.58'
$ws.Range("D45").Style = $cellStyle
$ws.Range("E45").Value = '  -3.50%  '
$ws.Range("E46").Value = '  -2.56%  '
$cellStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.92'
$ws.Range("D47").Style = $cellStyle
$ws.Range("E47").Value = '  -2.15%  '
$ws.Range("E48").Value = '  -3.85%  '
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("E50").Value = '  -0.11%  '
$cellStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.04'
$ws.Range("D51").Style = $cellStyle
$ws.Range("E51").Value = '  +0.31%  '
